# Applies the "fix eu tabelle for faktenblatt and other small things" commit
# to the faktenblatttabellen workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "R-Wert und 7-Tage-Inzidenz": drop the
#    "Regionen mit 7-TI bei Über-80-Jährigen:" block (rows 12-14).
# ---------------------------------------------------------------------------
$wsRWert = $wb.Worksheets.Item(2)
$wsRWert.Rows("12:14").Delete()

# ---------------------------------------------------------------------------
# 2) Sheet "Intensivbetten": drop the stray comma after the "%" figures.
# ---------------------------------------------------------------------------
$wsIntensiv = $wb.Worksheets.Item(3)
$wsIntensiv.Range("B3").Value = "15 %`n4108"
$wsIntensiv.Range("C3").Value = "17 %`n4552"
$wsIntensiv.Range("B4").Value = "20 %`n5539"
$wsIntensiv.Range("C4").Value = "19 %`n5198"

# ---------------------------------------------------------------------------
# 3) Sheet "Todesfälle und Fallsterblichkei": refreshed Übersterblichkeit
#    figures (rows 7-10).
# ---------------------------------------------------------------------------
$wsTod = $wb.Worksheets.Item(5)
$wsTod.Range("B7").Value = "-114 (-6,8%)"
$wsTod.Range("C7").Value = "-183 (-10,7%)"
$wsTod.Range("D7").Value = "60,5%"

$wsTod.Range("C8").Value = "-132 ( -2,1%)"
$wsTod.Range("D8").Value = "97,0%"

$wsTod.Range("B9").Value = "1218 (12,3%)"
$wsTod.Range("C9").Value = "1658 ( 16,7%)"
$wsTod.Range("D9").Value = "36,1%"

$wsTod.Range("C10").Value = "1344 (  7,5%)"

# ---------------------------------------------------------------------------
# 4) Sheet "Internationaler Vergleich": updated EU comparison table
#    (columns B "Anteil Bevölk." and C "Anzahl Fälle", rows 2-11).
# ---------------------------------------------------------------------------
$wsEU = $wb.Worksheets.Item(8)

$wsEU.Range("B2").Value = "5,4 %"
$wsEU.Range("C2").Value = 579079

$wsEU.Range("B3").Value = "5,3 %"
$wsEU.Range("C3").Value = 608137

$wsEU.Range("B4").Value = "3,7 %"
$wsEU.Range("C4").Value = 2432559

$wsEU.Range("B5").Value = "3,7 %"
$wsEU.Range("C5").Value = 1730575

$wsEU.Range("B6").Value = "3,6 %"
$wsEU.Range("C6").Value = 623567

$wsEU.Range("B7").Value = "3 %"
$wsEU.Range("C7").Value = 1843712

$wsEU.Range("B8").Value = "3 %"
$wsEU.Range("C8").Value = 1135676

$wsEU.Range("B9").Value = "2,9 %"
$wsEU.Range("C9").Value = 556335

$wsEU.Range("B10").Value = "2,7 %"
$wsEU.Range("C10").Value = 1854490

$wsEU.Range("B11").Value = "1,6 %"
$wsEU.Range("C11").Value = 1350810
